$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pass 1: write first occurrence of each NEW distinct string, in the exact order
# required to reproduce the target shared-string table ordering.
$ws.Cells.Item(29, 1).Value = 'PROJECTILE_VOLTAGE_FB'
$ws.Cells.Item(31, 1).Value = 'QSPI'
$ws.Cells.Item(29, 2).Value = 'PC06'
$ws.Cells.Item(30, 2).Value = 'PC05'
$ws.Cells.Item(30, 1).Value = 'PROJECTILE_CURRENT_FB'
$ws.Cells.Item(31, 2).Value = 'PB11'
$ws.Cells.Item(32, 2).Value = 'PB10'
$ws.Cells.Item(33, 2).Value = 'PA11'
$ws.Cells.Item(34, 2).Value = 'PA10'
$ws.Cells.Item(35, 2).Value = 'PA09'
$ws.Cells.Item(36, 2).Value = 'PA08'
$ws.Cells.Item(37, 2).Value = 'PA07'
$ws.Cells.Item(37, 1).Value = 'STATUS_LED_OP'
$ws.Cells.Item(38, 2).Value = 'PA06'
$ws.Cells.Item(38, 1).Value = 'ERROR_LED_OP'
$ws.Cells.Item(39, 2).Value = 'PA05'
$ws.Cells.Item(39, 3).Value = 'DAC'
$ws.Cells.Item(39, 4).Value = 'VOUT1'
$ws.Cells.Item(40, 4).Value = 'VOUT0'
$ws.Cells.Item(41, 3).Value = 'ADC IN'
$ws.Cells.Item(41, 4).Value = 'ADC1 1'
$ws.Cells.Item(41, 2).Value = 'PB09'
$ws.Cells.Item(41, 1).Value = 'PROJECTILE_CHARGE_CURRENT_AN'
$ws.Cells.Item(42, 1).Value = 'PROJECTILE_CHARGE_VOLTAGE_AN'
$ws.Cells.Item(42, 2).Value = 'PB08'
$ws.Cells.Item(42, 4).Value = 'ADC1 0'
$ws.Cells.Item(43, 1).Value = 'RES_PRESSURE_B_AN'
$ws.Cells.Item(43, 2).Value = 'PD00'
$ws.Cells.Item(43, 4).Value = 'ADC1 14'
$ws.Cells.Item(44, 1).Value = 'RES_PRESSURE_A_AN'
$ws.Cells.Item(44, 2).Value = 'PB05'
$ws.Cells.Item(44, 4).Value = 'ADC1 7'
$ws.Cells.Item(45, 1).Value = 'HW_VERSION_AN'
$ws.Cells.Item(45, 2).Value = 'PB04'
$ws.Cells.Item(45, 4).Value = 'ADC1 6'
$ws.Cells.Item(40, 2).Value = 'PA02'
$ws.Cells.Item(39, 1).Value = 'PROJECTILE_CHARGE_CURRENT_AO'
$ws.Cells.Item(40, 1).Value = 'PROJECTILE_CHARGE_VOLTAGE_AO'

# Pass 2: remaining cell writes (duplicate references to strings already
# present in the shared string table), in natural reading order.
$ws.Cells.Item(29, 3).Value = 'GPIO in, pull up.'
$ws.Cells.Item(30, 3).Value = 'GPIO in, pull up.'
$ws.Cells.Item(32, 1).Value = 'QSPI'
$ws.Cells.Item(33, 1).Value = 'QSPI'
$ws.Cells.Item(34, 1).Value = 'QSPI'
$ws.Cells.Item(35, 1).Value = 'QSPI'
$ws.Cells.Item(36, 1).Value = 'QSPI'
$ws.Cells.Item(37, 3).Value = 'GPIO out, default high.'
$ws.Cells.Item(38, 3).Value = 'GPIO out, default high.'
$ws.Cells.Item(40, 3).Value = 'DAC'
$ws.Cells.Item(42, 3).Value = 'ADC IN'
$ws.Cells.Item(43, 3).Value = 'ADC IN'
$ws.Cells.Item(44, 3).Value = 'ADC IN'
$ws.Cells.Item(45, 3).Value = 'ADC IN'

# Widen column A to fit new longer pin names (target OOXML width 32.7109375;
# the host's character->pixel rounding only lands on multiples of 1/6, so 31.9
# is the input that rounds to the closest achievable stored width, 32.667).
$ws.Columns.Item(1).ColumnWidth = 31.9

# Update view: scroll so row 16 is at top, and select A41 (mirrors author's final cursor position)
$ws.Range("A41").Select()
$excel.ActiveWindow.ScrollRow = 16
